$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 40: Number of Connected Components in an Undirected Graph ---
$ws.Cells.Item(40, 1).Value = "Graphs"

$ws.Cells.Item(40, 2).Value = "Medium"
$ws.Cells.Item(40, 2).Style = "Neutral"

$ws.Cells.Item(40, 3).Value = "Number of Connected Components in an Undirected Graph"
$ws.Cells.Item(40, 3).Style = "Neutral"
$ws.Hyperlinks.Add($ws.Cells.Item(40, 3), "https://neetcode.io/problems/count-connected-components?list=neetcode150")

$ws.Cells.Item(40, 4).Value = "We'll either use Union Find or DFS/BFS`nStart a DFS/BFS at every unvisited node; each start is one connected`ncomponent, and the DFS just marks all nodes in that component as seen`ncycles don’t matter.`nUnion-Find:`n- Initially, each node is its own parent (each node is its own component) and we also keep a `"rank`" array that stores the size of each component's root.`n- For every edge [a, b], we check if their roots are the same:`n        If yes → they're already in one component.`n        If no  → this edge connects two different components, so we merge them.`n- We merge smaller component into the bigger one (union by size) by attaching the smaller root to the bigger root and adding their sizes.`n- During find(), we apply path compression:`n        We recursively climb up until we find the root,`n        and along the way we set each node's parent directly to the root.`n    This flattens the structure, making all future finds extremely fast."
$ws.Cells.Item(40, 4).WrapText = $true
$ws.Cells.Item(40, 4).VerticalAlignment = -4160
$ws.Rows.Item(40).RowHeight = 43.8

# --- Row 41: 684. Redundant Connection ---
$ws.Cells.Item(41, 1).Value = "Graphs"

$ws.Cells.Item(41, 2).Value = "Medium"
$ws.Cells.Item(41, 2).Style = "Neutral"

$ws.Cells.Item(41, 3).Value = "684. Redundant Connection"
$ws.Cells.Item(41, 3).Style = "Neutral"
$ws.Hyperlinks.Add($ws.Cells.Item(41, 3), "https://leetcode.com/problems/redundant-connection/", "", "", "https://leetcode.com/problems/redundant-connection/")

$ws.Cells.Item(41, 4).Value = "We'll use Union Find.`nThe graph was initially a tree so each node WILL NOT have the same parent when we're building the tree otherwise there would be a cycle`nIf two nodes have the same parent then connecting them would cause a cycle, and it would`ndisqualify the tree property of the graph, so return the last one that causes a cycle"
$ws.Cells.Item(41, 4).WrapText = $true
$ws.Cells.Item(41, 4).VerticalAlignment = -4160
$ws.Rows.Item(41).RowHeight = 57.6

# --- Update selection to mirror the authored state (D42, one past the new last row) ---
$ws.Range("D42").Select()
